$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 27069.334
$ws.Range("I21").Value = 27069.334
$ws.Range("K21").Value = 27069.334
$ws.Range("M21").Value = -26601.334

# Row 23
$ws.Range("H23").Value = 27069.334
$ws.Range("I23").Value = 27069.334
$ws.Range("K23").Value = 27069.334
$ws.Range("M23").Value = -26835.334

# Row 31
$ws.Range("H31").Value = 119
$ws.Range("I31").Value = 133.25
$ws.Range("K31").Value = 399.75
$ws.Range("M31").Value = -169.75

# Row 33
$ws.Range("H33").Value = 117.416664
$ws.Range("J33").Value = 119.833336
$ws.Range("L33").Value = 119.833336
$ws.Range("N33").Value = -577.833336

# Row 34
$ws.Range("H34").Value = 2960.75
$ws.Range("I34").Value = 2960.75
$ws.Range("K34").Value = 2960.75
$ws.Range("M34").Value = -2757.75

# Row 36
$ws.Range("H36").Value = 2960.75
$ws.Range("I36").Value = 2960.75
$ws.Range("K36").Value = 2960.75
$ws.Range("M36").Value = -2245.75

# Row 45
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = ""
$ws.Range("N45").Value = -6384

# Row 100
$ws.Range("H100").Value = 2858.1667
$ws.Range("I100").Value = 2858.1667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2858.1667
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2317.1667
$ws.Range("N100").Value = ""

# Row 124
$ws.Range("H124").Value = 46442
$ws.Range("J124").Value = 46442
$ws.Range("L124").Value = 46442
$ws.Range("N124").Value = -56262

# Row 125
$ws.Range("H125").Value = 1428.8889
$ws.Range("J125").Value = 2800
$ws.Range("L125").Value = 25200
$ws.Range("N125").Value = -30120

# Row 132
$ws.Range("H132").Value = 1147.6
$ws.Range("I132").Value = 1147.6
$ws.Range("K132").Value = 3442.8
$ws.Range("M132").Value = -912.7999999999997

# Row 137
$ws.Range("H137").Value = 1428.3462
$ws.Range("I137").Value = 1223.3914
$ws.Range("K137").Value = 3670.1742
$ws.Range("M137").Value = -1120.1742

# Row 141
$ws.Range("H141").Value = 1871018.4
$ws.Range("I141").Value = 3503399
$ws.Range("K141").Value = 10510197
$ws.Range("M141").Value = -10505017

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4356.592
$ws.Range("I32").Value = 2385.6052
$ws.Range("J32").Value = 11165.454
$ws.Range("K32").Value = 2385.6052
$ws.Range("L32").Value = 11165.454
$ws.Range("M32").Value = -2098.6052
$ws.Range("N32").Value = -11739.454

# Row 45
$ws.Range("H45").Value = 1591.7142
$ws.Range("I45").Value = 1040.4
$ws.Range("K45").Value = 1040.4
$ws.Range("M45").Value = -663.4000000000001

# Row 61
$ws.Range("H61").Value = 3431.7407
$ws.Range("J61").Value = 11563.333
$ws.Range("L61").Value = 11563.333
$ws.Range("N61").Value = -11987.333

# Row 136
$ws.Range("H136").Value = 3431.7407
$ws.Range("J136").Value = 11563.333
$ws.Range("L136").Value = 34689.999
$ws.Range("N136").Value = -39789.999

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -887
$ws.Range("N5").Value = ""

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

# Row 132
$ws.Range("H132").Value = 45000
$ws.Range("J132").Value = 45000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -55120

# Row 134
$ws.Range("H134").Value = 7658.243
$ws.Range("I134").Value = 7262.1514
$ws.Range("J134").Value = 10926
$ws.Range("K134").Value = 21786.4542
$ws.Range("L134").Value = 32778
$ws.Range("M134").Value = -19251.4542
$ws.Range("N134").Value = -37848

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 104.1
$ws.Range("J7").Value = 68.333336
$ws.Range("L7").Value = 68.333336
$ws.Range("N7").Value = -294.333336

# Row 17
$ws.Range("H17").Value = 8019.8
$ws.Range("I17").Value = 99
$ws.Range("K17").Value = 99
$ws.Range("M17").Value = 75

# Row 29
$ws.Range("H29").Value = 7499.5
$ws.Range("J29").Value = 7499.5
$ws.Range("L29").Value = 7499.5
$ws.Range("N29").Value = -8085.5

# Row 31
$ws.Range("H31").Value = 2379.025
$ws.Range("I31").Value = 1601.5518
$ws.Range("K31").Value = 1601.5518
$ws.Range("M31").Value = -1306.5518

# Row 34
$ws.Range("H34").Value = 2379.025
$ws.Range("I34").Value = 1601.5518
$ws.Range("K34").Value = 1601.5518
$ws.Range("M34").Value = -1399.5518

# Row 132
$ws.Range("H132").Value = 2503.4
$ws.Range("I132").Value = 1402.3334
$ws.Range("K132").Value = 4207.0002
$ws.Range("M132").Value = -1677.0002

# Row 134
$ws.Range("H134").Value = 955.1429000000001
$ws.Range("I134").Value = 955.1429000000001
$ws.Range("K134").Value = 2865.4287
$ws.Range("M134").Value = -330.4287000000004

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 798.5
$ws.Range("I5").Value = 667.1667
$ws.Range("K5").Value = 2001.5001
$ws.Range("M5").Value = -1889.5001

# Row 11
$ws.Range("H11").Value = 659.5
$ws.Range("I11").Value = 383.33334
$ws.Range("K11").Value = 1150.00002
$ws.Range("M11").Value = -1010.00002

# Row 116
$ws.Range("H116").Value = 3082
$ws.Range("I116").Value = 1329
$ws.Range("J116").Value = 3666.3333
$ws.Range("K116").Value = 3987
$ws.Range("L116").Value = 10998.9999
$ws.Range("M116").Value = -545
$ws.Range("N116").Value = -17882.9999

# Row 117
$ws.Range("H117").Value = 509
$ws.Range("I117").Value = 509
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1527
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 1915
$ws.Range("N117").Value = ""

# Row 131
$ws.Range("H131").Value = 1394.19
$ws.Range("I131").Value = 589.875
$ws.Range("J131").Value = 1464.1305
$ws.Range("K131").Value = 1769.625
$ws.Range("L131").Value = 4392.3915
$ws.Range("M131").Value = 3270.375
$ws.Range("N131").Value = -14472.3915

# Row 135
$ws.Range("H135").Value = 798.5
$ws.Range("I135").Value = 667.1667
$ws.Range("K135").Value = 6004.5003
$ws.Range("M135").Value = -3469.5003

# Row 137
$ws.Range("H137").Value = 3469.7
$ws.Range("J137").Value = 4786.2
$ws.Range("L137").Value = 14358.6
$ws.Range("N137").Value = -24558.6

# Row 140
$ws.Range("H140").Value = 1916.4375
$ws.Range("I140").Value = 968.9
$ws.Range("J140").Value = 3495.6667
$ws.Range("K140").Value = 2906.7
$ws.Range("L140").Value = 10487.0001
$ws.Range("M140").Value = 2273.3
$ws.Range("N140").Value = -20847.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 34
$ws.Range("H34").Value = 20000
$ws.Range("J34").Value = 20000
$ws.Range("L34").Value = 20000
$ws.Range("N34").Value = -20536

# Row 76
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630

# Row 79
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184

# Row 132
$ws.Range("H132").Value = 2139758
$ws.Range("I132").Value = 2748878.5
$ws.Range("K132").Value = 8246635.5
$ws.Range("M132").Value = -8244105.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4733.2666
$ws.Range("I7").Value = 2267
$ws.Range("J7").Value = 8432.666999999999
$ws.Range("K7").Value = 2267
$ws.Range("L7").Value = 8432.666999999999
$ws.Range("M7").Value = -2155
$ws.Range("N7").Value = -8656.666999999999

# Row 43
$ws.Range("H43").Value = 12673.667
$ws.Range("J43").Value = 12673.667
$ws.Range("L43").Value = 12673.667
$ws.Range("N43").Value = -13059.667

# Row 93
$ws.Range("H93").Value = 596.3333
$ws.Range("I93").Value = 444.5
$ws.Range("K93").Value = 444.5
$ws.Range("M93").Value = 803.5

# Row 94
$ws.Range("H94").Value = 48000
$ws.Range("J94").Value = 48000
$ws.Range("L94").Value = 48000
$ws.Range("N94").Value = -49352

# Row 126
$ws.Range("H126").Value = 4733.2666
$ws.Range("I126").Value = 2267
$ws.Range("J126").Value = 8432.666999999999
$ws.Range("K126").Value = 6801
$ws.Range("L126").Value = 25298.001
$ws.Range("M126").Value = -4331
$ws.Range("N126").Value = -30238.001

# Row 132
$ws.Range("H132").Value = 2063.7917
$ws.Range("I132").Value = 1989.3
$ws.Range("K132").Value = 5967.9
$ws.Range("M132").Value = -3437.9

# Row 136
$ws.Range("H136").Value = 4119.1665
$ws.Range("I136").Value = 2536
$ws.Range("K136").Value = 7608
$ws.Range("M136").Value = -5058

# Row 138
$ws.Range("H138").Value = 67947.336
$ws.Range("J138").Value = 67947.336
$ws.Range("L138").Value = 67947.336
$ws.Range("N138").Value = -78227.336

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 126435.664
$ws.Range("I122").Value = 126435.664
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 379306.992
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -376856.992
$ws.Range("N122").Value = ""

# Row 123
$ws.Range("H123").Value = 47544.445
$ws.Range("J123").Value = 47544.445
$ws.Range("L123").Value = 47544.445
